$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current blank/summary block (row 112),
# shifting the blank separator row and the three summary rows down to 116-119.
$ws.Rows("112:115").Insert()

# The last working day (2014-04-01) gets a second entry split across two rows
# (C110/C111 day number 1 -> 2), and the evening shift on row 111 now runs
# later (ends 22:30 instead of 22:00).
$ws.Range("C110").Value = 2
$ws.Range("C111").Value = 2
$ws.Range("E111").Value = 0.9375

# New timesheet entries for 2014-04-03 and 2014-04-04 (PQ/PV-bus experiment).
$ws.Range("A112").Value = 2014
$ws.Range("B112").Value = 4
$ws.Range("C112").Value = 3
$ws.Range("D112").Value = 0.51736111111111105
$ws.Range("E112").Value = 0.53472222222222221
$ws.Range("F112").Formula = "=(E112-D112)*24*60"
$ws.Range("G112").Formula = "=F112/60"

$ws.Range("A113").Value = 2014
$ws.Range("B113").Value = 4
$ws.Range("C113").Value = 3
$ws.Range("D113").Value = 0.60416666666666663
$ws.Range("E113").Value = 0.64583333333333337
$ws.Range("F113").Formula = "=(E113-D113)*24*60"
$ws.Range("G113").Formula = "=F113/60"

$ws.Range("A114").Value = 2014
$ws.Range("B114").Value = 4
$ws.Range("C114").Value = 4
$ws.Range("D114").Value = 0.63888888888888895
$ws.Range("E114").Value = 0.75
$ws.Range("F114").Formula = "=(E114-D114)*24*60"
$ws.Range("G114").Formula = "=F114/60"

$ws.Range("A115").Value = 2014
$ws.Range("B115").Value = 4
$ws.Range("C115").Value = 4
$ws.Range("D115").Value = 0.84375
$ws.Range("E115").Value = 0.91666666666666663
$ws.Range("F115").Formula = "=(E115-D115)*24*60"
$ws.Range("G115").Formula = "=F115/60"

# Match the author's final selection.
[void]$ws.Range("E115").Select()
